$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell value updates from the source diff.
# Values are set directly (cell-by-cell) to exactly reproduce the target state:
# quantity (F) and recomputed value (G) changes, batch-row swaps (B/E/F/G),
# and the resulting Sub Total / Grand Total rollups (column B).

$ws.Range("F16").Value = 65
$ws.Range("G16").Value = 3329.95
$ws.Range("B19").Value = 3431.64
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("B52").Value = 6499.37
$ws.Range("F71").Value = 337
$ws.Range("G71").Value = 21466.9
$ws.Range("F86").Value = 70
$ws.Range("G86").Value = 8782.9
$ws.Range("B90").Value = 189341.73
$ws.Range("F115").Value = 214
$ws.Range("G115").Value = 20717.34
$ws.Range("B117").Value = 14672.42
$ws.Range("B127").Value = 57552
$ws.Range("E127").Value = 136.86
$ws.Range("F127").Value = -5
$ws.Range("G127").Value = -603.45
$ws.Range("B128").Value = 64329
$ws.Range("E128").Value = 128.32
$ws.Range("F128").Value = 2
$ws.Range("G128").Value = 241.38
$ws.Range("F145").Value = 562
$ws.Range("G145").Value = 4490.38
$ws.Range("B147").Value = 16732.96
$ws.Range("F163").Value = 12
$ws.Range("G163").Value = 3191.52
$ws.Range("F172").Value = 3
$ws.Range("G172").Value = 2016.12
$ws.Range("B175").Value = 30320.13
$ws.Range("B192").Value = 64973
$ws.Range("E192").Value = 35.4
$ws.Range("F192").Value = 2
$ws.Range("G192").Value = 66.59999999999999
$ws.Range("B193").Value = 48706
$ws.Range("E193").Value = 39.8
$ws.Range("F193").Value = -144
$ws.Range("G193").Value = -4795.2
$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 66
$ws.Range("G227").Value = 9522.48
$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32
$ws.Range("F231").Value = 3
$ws.Range("G231").Value = 355.38
$ws.Range("F247").Value = 153
$ws.Range("G247").Value = 15898.23
$ws.Range("F256").Value = 286
$ws.Range("G256").Value = 43234.62
$ws.Range("B260").Value = 197497.77
$ws.Range("F270").Value = 27
$ws.Range("G270").Value = 870.48
$ws.Range("B275").Value = 5900.64
$ws.Range("F288").Value = 45
$ws.Range("G288").Value = 4184.55
$ws.Range("F296").Value = 59
$ws.Range("G296").Value = 1250.8
$ws.Range("F302").Value = 58
$ws.Range("G302").Value = 12231.62
$ws.Range("B304").Value = 182426.3
$ws.Range("F320").Value = 58
$ws.Range("G320").Value = 3981.7
$ws.Range("F327").Value = 13
$ws.Range("G327").Value = 3269.5
$ws.Range("B330").Value = 28831.42
$ws.Range("B364").Value = 53602
$ws.Range("E364").Value = 15.69
$ws.Range("F364").Value = -231
$ws.Range("G364").Value = -3037.65
$ws.Range("B365").Value = 65068
$ws.Range("E365").Value = 13.97
$ws.Range("F365").Value = 63
$ws.Range("G365").Value = 828.45
$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9
$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29
$ws.Range("B375").Value = 45718
$ws.Range("E375").Value = 19.38
$ws.Range("F375").Value = -294
$ws.Range("G375").Value = -4768.68
$ws.Range("B376").Value = 64927
$ws.Range("E376").Value = 17.26
$ws.Range("F376").Value = 106
$ws.Range("G376").Value = 1719.32
$ws.Range("B382").Value = 45702
$ws.Range("E382").Value = 31.43
$ws.Range("F382").Value = -215
$ws.Range("G382").Value = -5654.5
$ws.Range("B383").Value = 64919
$ws.Range("E383").Value = 27.97
$ws.Range("F383").Value = 61
$ws.Range("G383").Value = 1604.3
$ws.Range("B385").Value = 65067
$ws.Range("E385").Value = 15.65
$ws.Range("F385").Value = 126
$ws.Range("G385").Value = 1855.98
$ws.Range("B386").Value = 53595
$ws.Range("E386").Value = 17.61
$ws.Range("F386").Value = -335
$ws.Range("G386").Value = -4934.55
$ws.Range("B442").Value = 64810
$ws.Range("E442").Value = 291.22
$ws.Range("F442").Value = 4
$ws.Range("G442").Value = 1095.68
$ws.Range("B443").Value = 53319
$ws.Range("E443").Value = 310.64
$ws.Range("F443").Value = -6
$ws.Range("G443").Value = -1643.52
$ws.Range("B473").Value = 60022
$ws.Range("E473").Value = 37.22
$ws.Range("F473").Value = -113
$ws.Range("G473").Value = -3709.79
$ws.Range("B474").Value = 64830
$ws.Range("E474").Value = 34.9
$ws.Range("F474").Value = 108
$ws.Range("G474").Value = 3545.64
$ws.Range("F485").Value = 19
$ws.Range("G485").Value = 3333.93
$ws.Range("B488").Value = 31271.93
$ws.Range("F509").Value = 233
$ws.Range("G509").Value = 18728.54
$ws.Range("B510").Value = 24861
$ws.Range("B572").Value = 65362
$ws.Range("F572").Value = 20
$ws.Range("G572").Value = 817.4
$ws.Range("B573").Value = 65079
$ws.Range("F573").Value = 6
$ws.Range("G573").Value = 245.22
$ws.Range("F578").Value = 90
$ws.Range("G578").Value = 4490.1
$ws.Range("F579").Value = 35
$ws.Range("G579").Value = 2821
$ws.Range("F582").Value = 41
$ws.Range("G582").Value = 2336.59
$ws.Range("B583").Value = 18591.2
$ws.Range("F592").Value = 0
$ws.Range("G592").Value = 0
$ws.Range("B593").Value = 5734.68
$ws.Range("F599").Value = 1798
$ws.Range("G599").Value = 293271.78
$ws.Range("F601").Value = 433
$ws.Range("G601").Value = 122482.71
$ws.Range("B606").Value = 465060.29
$ws.Range("B619").Value = 1838446.22
$ws.Range("B620").Value = 1838446.22
